$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels in row 3 ---
$ws.Range("C3").Value = "Downscaled"
$ws.Range("D3").Value = "Wavelet"
$ws.Range("E3").Value = "2D DCT"

# --- Clear old data area beyond the new extent ---
$ws.Range("B15:E18").ClearContents()

# --- Write the new data grid B4:E14 ---
$arr = New-Object 'object[,]' 11,4
$arr[0,0] = 784
$arr[0,1] = 0.97966
$arr[0,2] = 0.954833348
$arr[0,3] = 0.9683999
$arr[1,0] = 676
$arr[1,1] = $null
$arr[1,2] = 0.95856666
$arr[1,3] = 0.969066
$arr[2,0] = 576
$arr[2,1] = $null
$arr[2,2] = 0.9542333
$arr[2,3] = $null
$arr[3,0] = 484
$arr[3,1] = $null
$arr[3,2] = 0.959333
$arr[3,3] = 0.969125
$arr[4,0] = 400
$arr[4,1] = $null
$arr[4,2] = 0.95596665
$arr[4,3] = $null
$arr[5,0] = 324
$arr[5,1] = $null
$arr[5,2] = 0.9559666
$arr[5,3] = $null
$arr[6,0] = 256
$arr[6,1] = $null
$arr[6,2] = 0.950200001
$arr[6,3] = 0.970400005
$arr[7,0] = 196
$arr[7,1] = 0.96077495
$arr[7,2] = 0.95156667
$arr[7,3] = 0.9674500077
$arr[8,0] = 144
$arr[8,1] = $null
$arr[8,2] = 0.9543
$arr[8,3] = $null
$arr[9,0] = 100
$arr[9,1] = $null
$arr[9,2] = 0.9460333188
$arr[9,3] = 0.9651750028
$arr[10,0] = 64
$arr[10,1] = $null
$arr[10,2] = 0.923933
$arr[10,3] = $null
$ws.Range("B4:E14").Value = $arr

# --- Update the selection shown in the sheet view ---
$ws.Range("D13:D14").Select()
